# Updates the single-column results table to the refreshed DaCapo /
# Shenandoah GC / jython benchmark numbers.
#
# Rows 1-12 (1-based table row index) get new scalar values; the three
# multi-run "raw sample line" rows at the bottom of the table (previously
# holding tab-separated dumps) collapse down to the same three scalar
# values that used to live in rows 1-3.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "217"
    5  = "0.00002"
    6  = "0.00050"
    7  = "0.00017"
    9  = "0.00029"
    10 = "0.00034"
    11 = "0.00039"
    12 = "0.04323"
}

foreach ($rowIndex in $updates.Keys) {
    $t.Cell($rowIndex, 1).Range.Text = $updates[$rowIndex]
}

# The final three rows held tab-separated multi-column dumps; they are
# replaced with single scalar values (re-using the old rows 1-3 values).
$rowCount = $t.Rows.Count
$t.Cell($rowCount - 2, 1).Range.Text = "99.93"
$t.Cell($rowCount - 1, 1).Range.Text = "0.04"
$t.Cell($rowCount, 1).Range.Text = "64"
